# Generate Report for Handoff
# Updates the "Status" columns from "In Translation" to "Ready for handoff"
# and bumps the related handoff timestamps, across the Overview, zh-cn and
# de-de sheets of the localization-status workbook.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-03-23 16:40:20"

# zh-cn sheet: Status (C2) + Latest Handoff Datetime (E2)
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-23 16:40:11"

# de-de sheet: Status (C2) + Latest Handoff Datetime (E2)
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-23 16:40:20"
